$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 previously only had cells in columns B and E (mapping-file refs tied
# to the "contaminacion" / "ruidos-exteriores" columns). Those two columns
# move to D and G below, so fully remove the old B5/E5 cells (Clear, not
# ClearContents, so no empty placeholder cell is left behind).
$ws.Range("B5").Clear()
$ws.Range("E5").Clear()

# Row 1: prettified Spanish column headers (same 10 columns, reordered +
# relabelled: numero-hogares, contaminacion, comarca-codigo, municipio-nombre,
# ruidos-exteriores, provincia-nombre, aragon, provincia-codigo,
# comarca-nombre, municipio-codigo -> Comarca nombre, Numero hogares, ...).
$ws.Range("A1").Value = "Comarca nombre"
$ws.Range("B1").Value = "Número hogares"
$ws.Range("C1").Value = "Comarca código"
$ws.Range("D1").Value = "Contaminación"
$ws.Range("E1").Value = "Provincia código"
$ws.Range("F1").Value = "Aragón"
$ws.Range("G1").Value = "Ruidos exteriores"
$ws.Range("H1").Value = "Municipio código"
$ws.Range("I1").Value = "Provincia nombre"
$ws.Range("J1").Value = "Municipio nombre"

# Row 2: sdmx/iaest property reference per column, following the same
# reorder as row 1.
$ws.Range("A2").Value = "sdmx-dimension:refArea"
$ws.Range("B2").Value = "iaest-measure:numero-hogares"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "iaest-dimension:contaminacion"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("G2").Value = "iaest-dimension:ruidos-exteriores"
$ws.Range("H2").Value = "null"
$ws.Range("I2").Value = "sdmx-dimension:refArea"
$ws.Range("J2").Value = "sdmx-dimension:refArea"

# Row 3: dim / medida marker per column.
$ws.Range("A3").Value = "dim"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "dim"
$ws.Range("G3").Value = "dim"
$ws.Range("H3").Value = "null"
$ws.Range("I3").Value = "dim"
$ws.Range("J3").Value = "dim"

# Row 4: datatype / class / codelist URI per column.
$ws.Range("A4").Value = "URI-comarca"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "skos:Concept"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("G4").Value = "skos:Concept"
$ws.Range("H4").Value = "null"
$ws.Range("I4").Value = "URI-Provincia"
$ws.Range("J4").Value = "URI-Municipio"

# Row 5: mapping-file references, now under the Contaminación (D) and
# Ruidos exteriores (G) columns. Clone the standard cell style (s="1") from
# an existing data cell via copy/paste-special so no new style entry is
# introduced in styles.xml, then set the values.
$ws.Range("A1").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("G5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D5").Value = "mapping-contaminacion.xlsx"
$ws.Range("G5").Value = "mapping-ruidos-exteriores.xlsx"
